# Update org and User Roles
# Adds a generated-SQL helper column (AA) to the "result" sheet:
#   - AA1: header/label holding the INSERT INTO ... VALUES preamble string
#   - AA2: one-off formula building the SQL values tuple for row 2
#   - AA3:AA12: the same formula (relative refs) filled down, stored by
#     Excel as a shared formula group
# Also fixes the "level" value for the Incident Primary Contact role (V5)
# from 2 to 3, and moves the sheet selection/scroll to reflect where the
# author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label for the new helper column -------------------------------
$ws.Range("AA1").Value = "INSERT INTO public.user_roles(id, name_t, description_t, is_active, create_incident, approve_orgs_full, approve_orgs_preliminary, edit_portal_settings, affiliate_org, invite_workers, remove_workers, phone_agent, advanced_maps, translate, support_agent, crew_management, view_user_contacts, view_sensitive, is_default, level, created_at, updated_at) VALUES"

# --- Row 2: standalone (non-shared) formula --------------------------------
$ws.Range("AA2").Formula = "=""(""&A2&"", '""&D2&""', '""&E2&""', ""&IF(F2,""true"",""false"")&"", ""&IF(G2,""true"",""false"")&"", ""&IF(H2,""true"",""false"")&"", ""&IF(I2,""true"",""false"")&"", ""&IF(J2,""true"",""false"")&"", ""&IF(K2,""true"",""false"")&"", ""&IF(L2,""true"",""false"")&"", ""&IF(M2,""true"",""false"")&"", ""&IF(N2,""true"",""false"")&"", ""&IF(O2,""true"",""false"")&"", ""&IF(P2,""true"",""false"")&"", ""&IF(Q2,""true"",""false"")&"", ""&IF(R2,""true"",""false"")&"", ""&IF(S2,""true"",""false"")&"", '""&T2&""', ""&IF(U2,""true"",""false"")&"", ""&V2&"", '""&TEXT(C2,""YYYY-MM-DD"")&""', NOW()),"""

# --- Rows 3-12: same formula, filled down (becomes a shared formula group) -
$ws.Range("AA3:AA12").Formula = "=""(""&A3&"", '""&D3&""', '""&E3&""', ""&IF(F3,""true"",""false"")&"", ""&IF(G3,""true"",""false"")&"", ""&IF(H3,""true"",""false"")&"", ""&IF(I3,""true"",""false"")&"", ""&IF(J3,""true"",""false"")&"", ""&IF(K3,""true"",""false"")&"", ""&IF(L3,""true"",""false"")&"", ""&IF(M3,""true"",""false"")&"", ""&IF(N3,""true"",""false"")&"", ""&IF(O3,""true"",""false"")&"", ""&IF(P3,""true"",""false"")&"", ""&IF(Q3,""true"",""false"")&"", ""&IF(R3,""true"",""false"")&"", ""&IF(S3,""true"",""false"")&"", '""&T3&""', ""&IF(U3,""true"",""false"")&"", ""&V3&"", '""&TEXT(C3,""YYYY-MM-DD"")&""', NOW()),"""

# --- Fix the "level" value for the Incident Primary Contact role ----------
$ws.Range("V5").Value = 3

# --- Reflect the author's final selection/scroll position -----------------
[void]$ws.Range("AB10").Select()

$wb.Save()
